$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-13 Thursday" "2024-06-14 Friday"

Replace-Text "97÷4=24, 1" "13÷8=1, 5"
Replace-Text "21÷3=7, 0" "61÷3=20, 1"
Replace-Text "93÷3=31, 0" "85÷3=28, 1"
Replace-Text "17÷3=5, 2" "54÷8=6, 6"
Replace-Text "27÷6=4, 3" "95÷9=10, 5"

Replace-Text "99÷8=12, 3" "21÷2=10, 1"
Replace-Text "26÷4=6, 2" "14÷5=2, 4"
Replace-Text "92÷5=18, 2" "18÷3=6, 0"
Replace-Text "37÷5=7, 2" "43÷9=4, 7"
Replace-Text "41÷7=5, 6" "62÷6=10, 2"

Replace-Text "77÷4=19, 1" "41÷7=5, 6"
Replace-Text "89÷9=9, 8" "43÷2=21, 1"
Replace-Text "78÷9=8, 6" "52÷7=7, 3"
Replace-Text "44÷9=4, 8" "18÷7=2, 4"
Replace-Text "72÷5=14, 2" "92÷2=46, 0"

Replace-Text "95÷5=19, 0" "11÷7=1, 4"
Replace-Text "16÷6=2, 4" "61÷9=6, 7"
Replace-Text "78÷3=26, 0" "85÷2=42, 1"
Replace-Text "92÷3=30, 2" "14÷6=2, 2"
Replace-Text "73÷7=10, 3" "82÷3=27, 1"

Replace-Text "29÷6=4, 5" "87÷3=29, 0"
Replace-Text "25÷9=2, 7" "71÷5=14, 1"
Replace-Text "50÷8=6, 2" "42÷2=21, 0"
Replace-Text "80÷6=13, 2" "18÷2=9, 0"
Replace-Text "21÷4=5, 1" "11÷2=5, 1"

Write-Host "Replacements complete"
